$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated values
$ws.Range("B2").Value = 0.6606524410359556
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 3.56341032713086

# Row 3 updated values
$ws.Range("B3").Value = 0.003208871385164791
$ws.Range("C3").Value = 0.306821227259698
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 1.557009902479497
